$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44797
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16500
$ws.Range("P3").Value = 917

# Row 4
$ws.Range("D4").Value = 44809
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("P4").Value = 806

# Row 5
$ws.Range("D5").Value = 44659
$ws.Range("J5").Value = 80

# Row 6
$ws.Range("D6").Value = 44818
$ws.Range("J6").Value = 60
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("P6").Value = 833

# Row 7
$ws.Range("D7").Value = 44656
$ws.Range("J7").Value = 100

# Row 8
$ws.Range("D8").Value = 44799

# Row 9
$ws.Range("D9").Value = 44664
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 861

# Row 11
$ws.Range("D11").Value = 44658
$ws.Range("J11").Value = 80

# Row 12
$ws.Range("D12").Value = 44637

# Row 13
$ws.Range("D13").Value = 44785

# Row 14
$ws.Range("D14").Value = 44628
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15500
$ws.Range("P14").Value = 861

# Row 15
$ws.Range("D15").Value = 44384
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17500
$ws.Range("P15").Value = 972

# Row 16
$ws.Range("D16").Value = 44384
$ws.Range("I16").Value = "Segunda"
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("P16").Value = 833

# Row 18
$ws.Range("D18").Value = 44642
$ws.Range("J18").Value = 100

# Row 19
$ws.Range("D19").Value = 44804
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15500
$ws.Range("P19").Value = 861

# Row 20
$ws.Range("D20").Value = 44761
$ws.Range("J20").Value = 100

# Row 21
$ws.Range("D21").Value = 44819

# Row 22
$ws.Range("D22").Value = 44813
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 14000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 14500
$ws.Range("P22").Value = 806

# Row 23
$ws.Range("D23").Value = 44790
$ws.Range("I23").Value = "Primera"
$ws.Range("K23").Value = 17000
$ws.Range("L23").Value = 18000
$ws.Range("M23").Value = 17500
$ws.Range("P23").Value = 972

# Row 24
$ws.Range("D24").Value = 44771
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 17000
$ws.Range("L24").Value = 18000
$ws.Range("M24").Value = 17500
$ws.Range("P24").Value = 972

# Row 25
$ws.Range("D25").Value = 44651
$ws.Range("J25").Value = 60

# Row 26
$ws.Range("D26").Value = 44791
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17500
$ws.Range("P26").Value = 972

# Row 27
$ws.Range("D27").Value = 44811
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = 14500
$ws.Range("P27").Value = 806

# Row 28
$ws.Range("D28").Value = 44645
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 16000
$ws.Range("M28").Value = 15500
$ws.Range("P28").Value = 861

# Row 29
$ws.Range("D29").Value = 44630
$ws.Range("K29").Value = 15000
$ws.Range("L29").Value = 16000
$ws.Range("M29").Value = 15500
$ws.Range("P29").Value = 861

# Row 30
$ws.Range("D30").Value = 44635
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15500
$ws.Range("P30").Value = 861

# Row 31
$ws.Range("D31").Value = 44782
$ws.Range("J31").Value = 120

# Row 32
$ws.Range("D32").Value = 44763
$ws.Range("J32").Value = 80

# Row 33
$ws.Range("D33").Value = 44775
$ws.Range("J33").Value = 100
